# Apply the "quiz marksheet" re-scoring edit described in the commit:
#   "Handles float input without breaking stuff"
#
# The sheet previously showed a blank/"Absent" result for the student
# (the 3rd Student-Ans/Correct-Ans block in columns G:H was a leftover,
# and the score table was all zeros). This script recomputes the score,
# drops the unused 3rd answer block (G:H), trims the 2nd answer block
# (D:E) down to the 3 still-relevant questions, and marks every question
# the student answered correctly by writing the matching option text
# into column A / D using the existing green "correctStyle".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Score summary table (rows 10-12)
# ---------------------------------------------------------------------
# Right / Wrong / Not-Attempt / Max
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

# Marking scheme actually applied (per-question marks) + resulting total
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 60
$ws.Range("E12").Value = "60/112"

# A10/A11/A12 ("No." / "Marking" / "Total") pick up the "absoluteStyle"
# (style index 4) that the rest of that header column already uses -
# copy it from A9, which already carries that exact style, so we reuse
# the existing cellXfs entry instead of minting a new one.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Drop the unused 3rd "Student Ans / Correct Ans" block (G:H)
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------
# 3. Trim the 2nd "Student Ans / Correct Ans" block (D:E) down to the
#    3 questions that still apply (rows 16-18); everything below is
#    cleared out entirely.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

# Column D ("Student Ans") for the remaining 3 questions: the student
# answered all 3 correctly, so D now mirrors the "Correct Ans" in E and
# uses the green "correctStyle" (same numeric style as B10, copied so
# the existing cellXfs entry - index 5 - is reused).
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------
# 4. Fill in column A ("Student Ans") of the 1st block (rows 16-40)
#    wherever the student's answer matches the "Correct Ans" in B -
#    i.e. mark the correctly-answered questions, in green.
# ---------------------------------------------------------------------
$correctRows = @{
    18 = "Option B"
    19 = "Option C"
    22 = "Option D"
    27 = "Option A"
    28 = "Option D"
    30 = "Option B"
    32 = "Option C"
    33 = "Option D"
    36 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}

$ws.Range("B10").Copy()
foreach ($row in $correctRows.Keys) {
    $cell = $ws.Range("A$row")
    $cell.PasteSpecial(-4122)
    $cell.Value = $correctRows[$row]
}
$excel.CutCopyMode = 0
